# Apply the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Sun Apr 16 13:59:03 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '30.564.44'
$ws.Range("E2").Value = '  +0.34%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.105.71'
$ws.Range("E3").Value = '  -0.02%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.018'
$ws.Range("E4").Value = '  +1.35%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '346.12'
$ws.Range("E5").Value = '  +3.75%  '

# Row 6: USDC
$ws.Range("E6").Value = '  +0.99%  '

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5249'
$ws.Range("E7").Value = '  +0.15%  '

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4493'
$ws.Range("E8").Value = '  -2.02%  '

# Row 9: OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.65'

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08985'
$ws.Range("E10").Value = '  +0.27%  '

# Row 11: Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.166'
$ws.Range("E11").Value = '  -0.94%  '

# Row 12: Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.27'
$ws.Range("E12").Value = '  -0.34%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '2.115.23'
$ws.Range("E13").Value = '  +0.86%  '

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.787'
$ws.Range("E14").Value = '  +0.22%  '

# Row 15: Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.043'
$ws.Range("E15").Value = '  +2.64%  '

# Row 16: Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.35'
$ws.Range("E16").Value = '  +2.98%  '

# Row 17: ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001171'
$ws.Range("E17").Value = '  +3.86%  '

# Row 18: BinanceUSD
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.016'
$ws.Range("E18").Value = '  +1.10%  '

# Row 19: TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06718'
$ws.Range("E19").Value = '  +1.41%  '

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.27'
$ws.Range("E20").Value = '  +0.22%  '

# Row 21: Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("E21").Value = '  +1.04%  '

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.302'
$ws.Range("E22").Value = '  +0.30%  '

# Row 23: WrappedBTC
$ws.Range("D23").Value = '30.642.21'
$ws.Range("E23").Value = '  +0.37%  '

# Row 24: Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.71'
$ws.Range("E24").Value = '  +3.25%  '

# Row 25: Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.397'
$ws.Range("E25").Value = '  +1.46%  '

# Row 26: WrappedliquidstakedEther2.0
$ws.Range("D26").Value = '2.361.80'
$ws.Range("E26").Value = '  +0.66%  '

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.28'
$ws.Range("E27").Value = '  -0.08%  '

# Row 28: Monero
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '165.10'
$ws.Range("E28").Value = '  +1.12%  '

# Row 29: LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.507'
$ws.Range("E29").Value = '  -1.86%  '

# Row 30: BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.75'
$ws.Range("E30").Value = '  +2.43%  '

# Row 31: ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.183'
$ws.Range("E31").Value = '  -0.75%  '

# Row 32: Stellar
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1071'
$ws.Range("E32").Value = '  -0.04%  '

# Row 33: Filecoin
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.312'
$ws.Range("E33").Value = '  +2.81%  '

# Row 34: ARBITRUM
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.612'
$ws.Range("E34").Value = '  -4.19%  '

# Row 35: HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.006'
$ws.Range("E35").Value = '  +1.92%  '

# Row 36: FraxShare
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.18'
$ws.Range("E36").Value = '  -2.31%  '

# Row 37: InternetComputer(DFINITY)
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.868'
$ws.Range("E37").Value = '  +6.14%  '

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02630'
$ws.Range("E38").Value = '  +2.40%  '

# Row 39: Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06800'
$ws.Range("E39").Value = '  -0.13%  '

# Row 40: Algorand
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2304'
$ws.Range("E40").Value = '  +0.73%  '

# Row 41: Aptos
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.52'
$ws.Range("E41").Value = '  -1.68%  '

# Row 42: TheSandbox
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6829'
$ws.Range("E42").Value = '  -0.72%  '

# Row 43: TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.266'
$ws.Range("E43").Value = '  +1.78%  '

# Row 44: EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.80'
$ws.Range("E44").Value = '  +6.52%  '

# Row 45: Decentraland
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6391'
$ws.Range("E45").Value = '  +0.28%  '

# Row 46: NEARProtocol
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.305'
$ws.Range("E46").Value = '  -1.93%  '

# Row 47: PancakeSwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.745'
$ws.Range("E47").Value = '  +2.61%  '

# Row 48: BabyDogeCoin
$ws.Range("E48").Value = '  +2.57%  '

# Row 49: EOS
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.249'
$ws.Range("E49").Value = '  +0.48%  '

# Row 50: Cronos
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07290'
$ws.Range("E50").Value = '  +2.64%  '

# Row 51: WEMIXTOKEN
$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.194'
$ws.Range("E51").Value = '  -2.10%  '
